$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D numeric-looking strings are written as text (preserve exact formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.757.91"

$ws.Range("D3").Value = "1.547.48"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "206.14"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("E6").Value = "  -2.03%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "21.41"
$ws.Range("E8").Value = "  -3.75%  "

$ws.Range("E10").Value = "  -1.25%  "

$ws.Range("E11").Value = "  -1.71%  "

$ws.Range("D12").Value = "1.769.61"
$ws.Range("E12").Value = "  -1.58%  "

$ws.Range("D13").Value = "1.559.93"
$ws.Range("E13").Value = "  -0.99%  "

$ws.Range("E14").Value = "  -2.78%  "

$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  -1.28%  "

$ws.Range("D16").Value = "26.739.28"
$ws.Range("E16").Value = "  -1.60%  "

$ws.Range("D17").Value = "61.14"
$ws.Range("E17").Value = "  -1.76%  "

$ws.Range("D18").Value = "212.51"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("D23").Value = "8.93"
$ws.Range("E23").Value = "  -5.77%  "

$ws.Range("E24").Value = "  -1.44%  "

$ws.Range("D25").Value = "152.96"
$ws.Range("E25").Value = "  +0.36%  "

$ws.Range("D26").Value = "6.50"
$ws.Range("E26").Value = "  -3.19%  "

$ws.Range("D27").Value = "14.87"
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  -1.95%  "

$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("E31").Value = "  -1.55%  "

$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("D33").Value = "1.340.95"
$ws.Range("E33").Value = "  -4.09%  "

$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("E35").Value = "  -3.34%  "

$ws.Range("E36").Value = "  -0.58%  "

$ws.Range("D37").Value = "0.926"
$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("D39").Value = "0.520"
$ws.Range("E39").Value = "  +0.71%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "5.81"
$ws.Range("E40").Value = "  +7.08%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "0.800"
$ws.Range("E41").Value = "  -1.87%  "

$ws.Range("D42").Value = "0.994"
$ws.Range("E42").Value = "  -1.61%  "

$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("E44").Value = "  -4.20%  "

$ws.Range("D45").Value = "62.76"
$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("D46").Value = "1.682.42"
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("E47").Value = "  -3.94%  "

$ws.Range("D48").Value = "85.95"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").Value = "0.0509"
$ws.Range("E49").Value = "  +3.06%  "

$ws.Range("D50").Value = "0.0₇0979"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("D51").Value = "0.0952"
$ws.Range("E51").Value = "  -0.02%  "

# Restore default (unstyled) style for column D after forcing text format
$ws.Range("D2:D51").Style = "Normal"
